# Apply corrected unit_sales and price values for the affected weeks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;  UnitSales = 13460; Price = 4.99 },
    @{ Row = 8;  UnitSales = 13397; Price = 4.99 },
    @{ Row = 10; UnitSales = 13388; Price = 4.99 },
    @{ Row = 11; UnitSales = 13324; Price = 4.99 },
    @{ Row = 13; UnitSales = 13467; Price = 4.99 },
    @{ Row = 20; UnitSales = 13500; Price = 4.99 },
    @{ Row = 23; UnitSales = 12289; Price = 4.5  },
    @{ Row = 27; UnitSales = 13288; Price = 4.99 },
    @{ Row = 28; UnitSales = 13512; Price = 4.99 },
    @{ Row = 29; UnitSales = 13240; Price = 4.99 },
    @{ Row = 30; UnitSales = 12220; Price = 4.5  },
    @{ Row = 32; UnitSales = 16080; Price = 4.5  },
    @{ Row = 33; UnitSales = 12316; Price = 4.5  },
    @{ Row = 35; UnitSales = 12444; Price = 4.5  },
    @{ Row = 36; UnitSales = 16326; Price = 4.5  },
    @{ Row = 37; UnitSales = 13160; Price = 4.99 },
    @{ Row = 40; UnitSales = 13365; Price = 4.99 },
    @{ Row = 41; UnitSales = 12544; Price = 4.5  },
    @{ Row = 42; UnitSales = 13414; Price = 4.99 },
    @{ Row = 43; UnitSales = 13266; Price = 4.99 },
    @{ Row = 46; UnitSales = 12650; Price = 4.5  },
    @{ Row = 48; UnitSales = 16049; Price = 4.5  },
    @{ Row = 49; UnitSales = 12721; Price = 4.5  },
    @{ Row = 51; UnitSales = 13314; Price = 4.99 },
    @{ Row = 53; UnitSales = 12473; Price = 4.5  }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.UnitSales
    $ws.Cells.Item($u.Row, 6).Value = $u.Price
}
